$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$oldId = "e506dd5f-9a7d-49c3-a340-16088b0780e4"
$newId = "1958bf87-a8aa-4e62-96af-226375b5b202"

$oldHash = "a2570090d83950d7041113a6a216aace429495b6"
$newHash = "969ad0aaa490cffa3117e57ca16025d8cdb15281"

$newFileName = "$newId.md"
$newDisplay = "e2e\$newId.md"

$newGenerateDate = "2016-08-25 00:54:38"
$newZhCnDate = "2016-08-25 00:54:33"

$newZhCnXlf = "$newId.$newHash.zh-cn.xlf"
$newDeDeXlf = "$newId.$newHash.de-de.xlf"

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newDisplay
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newDisplay
$wsOverview.Range("G2").Value = $newGenerateDate

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = $newFileName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnDate

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = $newFileName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newGenerateDate
